$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F3").Value = 772
    $ws.Range("F6").Value = 4553
    $ws.Range("F8").Value = 369
    $ws.Range("F9").Value = 1316
    $ws.Range("F14").Value = 509
}
